$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.364.27'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.604.88'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''212.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''0.499'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.244'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = '''0.0605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '''19.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '1.832.33'
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("D13").Value = '1.599.67'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '''0.506'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '''63.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.87%  '
$ws.Range("D17").Value = '26.377.17'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").Value = '''230.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.42%  '
$ws.Range("D19").Value = '''7.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.61%  '
$ws.Range("D20").Value = '0.0₃0724'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D22").Value = '''4.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = '''8.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.01%  '
$ws.Range("D24").Value = '''2.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").Value = '''146.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '''6.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").Value = '''15.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").Value = '''0.0494'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = '1.485.13'
$ws.Range("E32").Value = '  +5.58%  '
$ws.Range("D33").Value = '''3.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").Value = '''2.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").Value = '''1.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '''0.560'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").Value = '''0.821'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("D43").Value = '''0.930'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.50%  '
$ws.Range("D44").Value = '1.743.30'
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").Value = '''60.83'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '''89.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''0.0960'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '''1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.07%  '
